# Applies the per-cell text updates captured in the authoritative diff for cryptos.xlsx.
# The sheet stores every data value as text (t="inlineStr"/shared string), including
# numeric-looking Price figures (column D) such as "39.50" or "0.0785". Plain
# `Range.Value = "..."` assignment on such strings gets auto-coerced to a real number by
# the Excel object model (dropping trailing zeros, flipping to scientific notation, etc.),
# so column-D cells are temporarily forced to Text format while the value is written, then
# restored to the default "Normal" style so no stray formatting is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = [ordered]@{
    "D2" = "47.329.60"
    "E2" = "  +0.42%  "
    "D3" = "2.493.31"
    "E3" = "  -0.05%  "
    "E4" = "  +0.09%  "
    "D5" = "321.25"
    "E5" = "  -0.50%  "
    "D6" = "108.62"
    "E6" = "  +3.27%  "
    "D7" = "0.523"
    "E7" = "  +0.17%  "
    "E8" = "  -0.01%  "
    "D9" = "0.537"
    "E9" = "  -0.98%  "
    "D10" = "39.50"
    "E10" = "  +5.49%  "
    "E11" = "  -0.39%  "
    "E12" = "  +0.15%  "
    "D13" = "18.39"
    "E13" = "  +0.57%  "
    "E14" = "  -1.29%  "
    "D15" = "2.883.03"
    "E15" = "  +0.08%  "
    "D16" = "2.502.39"
    "E16" = "  +0.43%  "
    "E17" = "  +0.40%  "
    "D18" = "47.250.94"
    "E18" = "  +0.45%  "
    "D19" = "13.20"
    "E19" = "  +4.54%  "
    "E20" = "  +1.06%  "
    "E21" = "  +0.13%  "
    "E22" = "  +12.02%  "
    "D23" = "70.38"
    "E23" = "  -0.80%  "
    "D24" = "245.16"
    "E25" = "  +0.82%  "
    "D27" = "25.73"
    "E27" = "  -1.73%  "
    "E28" = "  +0.42%  "
    "E29" = "  -1.52%  "
    "B30" = "InjectiveProtocol"
    "C30" = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
    "D30" = "34.83"
    "E30" = "  -1.44%  "
    "B31" = "Kaspa"
    "C31" = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
    "D31" = "0.137"
    "E31" = "  +1.74%  "
    "D33" = "20.66"
    "E33" = "  +4.81%  "
    "E34" = "  -0.11%  "
    "D35" = "0.0785"
    "E35" = "  +0.06%  "
    "E36" = "  +0.20%  "
    "B37" = "ARBITRUM"
    "C37" = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
    "D37" = "1.97"
    "E37" = "  +2.01%  "
    "B38" = "RenderToken"
    "C38" = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
    "D38" = "4.71"
    "E38" = "  +2.11%  "
    "D39" = "2.92"
    "E39" = "  -1.64%  "
    "D40" = "22.99"
    "E40" = "  +5.94%  "
    "E41" = "  -0.02%  "
    "E42" = "  +0.69%  "
    "D43" = "116.86"
    "E43" = "  -4.08%  "
    "D44" = "0.0297"
    "E44" = "  +0.65%  "
    "D45" = "1.998.21"
    "E45" = "  +2.35%  "
    "E46" = "  +2.05%  "
    "E47" = "  -4.83%  "
    "E48" = "  +0.49%  "
    "D49" = "1.77"
    "E49" = "  -0.81%  "
    "E50" = "  -4.80%  "
    "D51" = "56.47"
    "E51" = "  +3.38%  "
}

# Column-D (Price) cells: force Text format so numeric-looking strings keep their exact
# original representation instead of being parsed into a Double.
$priceCells = @(
    "D2",
    "D3",
    "D5",
    "D6",
    "D7",
    "D9",
    "D10",
    "D13",
    "D15",
    "D16",
    "D18",
    "D19",
    "D23",
    "D24",
    "D27",
    "D30",
    "D31",
    "D33",
    "D35",
    "D37",
    "D38",
    "D39",
    "D40",
    "D43",
    "D44",
    "D45",
    "D49",
    "D51",
)
foreach ($cellRef in $priceCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}

# Drop the temporary Text format again so the cell style matches the rest of the sheet
# (un-styled data cells, i.e. no explicit `s` attribute).
foreach ($cellRef in $priceCells) {
    $ws.Range($cellRef).Style = "Normal"
}

Write-Host "Applied $($updates.Count) cell updates"
